$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 158-160: Masked UIN card template (RPR_MASKED_UIN_CARD_TEMPLATE)
# Columns: A=code, B=descr, C=lang_code, D=is_active, E=cr_by, F=cr_dtimes
# Shared-string table order matches: descr(eng), descr(ara), descr(fra), code
# -> populate column B for all three rows before column A.

$ws.Range("B158").Value = "Masked UIN card template"
$ws.Range("B159").Value = "قالب بطاقة UIN Masked"
$ws.Range("B160").Value = "Modèle de carte Masked UIN"

$ws.Range("A158").Value = "RPR_MASKED_UIN_CARD_TEMPLATE"
$ws.Range("A159").Value = "RPR_MASKED_UIN_CARD_TEMPLATE"
$ws.Range("A160").Value = "RPR_MASKED_UIN_CARD_TEMPLATE"

$ws.Range("C158").Value = "eng"
$ws.Range("D158").Value = $true
$ws.Range("E158").Value = "superadmin"
$ws.Range("F158").Value = "now()"

$ws.Range("C159").Value = "ara"
$ws.Range("D159").Value = $true
$ws.Range("E159").Value = "superadmin"
$ws.Range("F159").Value = "now()"

$ws.Range("C160").Value = "fra"
$ws.Range("D160").Value = $true
$ws.Range("E160").Value = "superadmin"
$ws.Range("F160").Value = "now()"

# Apply alignment formatting to the code/descr columns of the new rows
# (mirrors the new cellXf with applyAlignment introduced in the source edit)
$ws.Range("A158:B160").IndentLevel = 0

# Update view state to match the post-edit selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 145
[void]$ws.Range("A160").Select()
